$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("First", "Second")
    ,@("MC", "PAC")
    ,@("MC", "SLC")
    ,@("DC", "MC")
    ,@("E3", "E5")
    ,@("RCH", "DWE")
    ,@("E2", "DWE")
    ,@("E3", "DC")
    ,@("C2", "DC")
    ,@("B1", "B2")
    ,@("B1", "ESC")
    ,@("B2", "QNC")
    ,@("C2", "C1")
    ,@("C2", "ESC")
    ,@("DC", "M3")
    ,@("E2", "E3")
    ,@("E2", "PHYS")
    ,@("ESC", "C2")
    ,@("EV1", "EV2")
    ,@("EV1", "ML")
    ,@("EV1", "AL")
    ,@("EV1", "TC")
    ,@("EV2", "EV3")
    ,@("M3", "MC")
    ,@("PAS", "EV2")
    ,@("PHYS", "E2")
    ,@("QNC", "B2")
    ,@("QNC", "SLC")
    ,@("QNC", "MC")
    ,@("RCH", "E2")
    ,@("SCH", "AL")
    ,@("SCH", "HH")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("A16:C16").Select() | Out-Null
